$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MAE)
$ws.Range("B2").Value = 0.723
$ws.Range("C2").Value = 0.592
$ws.Range("D2").Value = 0.601
$ws.Range("E2").Value = 1.43
$ws.Range("F2").Value = 0.866

# Row 3 (MSE)
$ws.Range("B3").Value = 1.138
$ws.Range("C3").Value = 0.672
$ws.Range("D3").Value = 0.698
$ws.Range("E3").Value = 7.504
$ws.Range("F3").Value = 1.69

# Row 4 (mean Y-Test)
$ws.Range("B4").Value = 18.214
$ws.Range("C4").Value = 15.308
$ws.Range("D4").Value = 12.948
$ws.Range("E4").Value = 30.588
$ws.Range("F4").Value = 18.064

# Row 5 (mean Y-predicted)
$ws.Range("B5").Value = 18.282
$ws.Range("C5").Value = 15.294
$ws.Range("D5").Value = 12.857
$ws.Range("E5").Value = 30.256
$ws.Range("F5").Value = 17.659

# Row 6 (R2)
$ws.Range("B6").Value = 0.903
$ws.Range("C6").Value = 0.96
$ws.Range("D6").Value = 0.871
$ws.Range("E6").Value = 0.827
$ws.Range("F6").Value = 0.924
